$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = 10613.61331315967
$ws.Range("D3").Value = 697.1794862465113

$ws.Range("B4").Value = 4525.229245257393
$ws.Range("D4").Value = 639.0346550209613

$ws.Range("B5").Value = 1069.025

$ws.Range("B6").Value = 7080.953

$ws.Range("B7").Value = 9676.126000000007
$ws.Range("D7").Value = 920

$ws.Range("B8").Value = 14378.45375000002
$ws.Range("D8").Value = 1280

$ws.Range("B9").Value = 24293.37300000006
$ws.Range("D9").Value = 5100

$ws.Range("F10").Value = 5578515157.397015

$ws.Range("G11").Value = 0.7210308521362605

$ws.Range("F12").Value = 366437537.9720001
$ws.Range("G12").Value = 0.0656872891142207

$ws.Range("G13").Value = 0.2132818587495187
